# DPLKINV119-001 - Jenis Porto - Reksadana Ubah Data
# Update the "Kode Jenis Porto" / "Nama Jenis Porto" test-data pair from
# R07 / Reksadana Campuran Syariah to R91 / Reksadana Konvensional, and
# move the active selection/scroll position over to the M:N columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# M2 = KODE_JENIS_PORTO, N2 = NAMA_JENIS_PORTO
$ws.Range("M2").Value = "R91"
$ws.Range("N2").Value = "Reksadana Konvensional"

# Move the view / selection so that column G is the left-most visible
# column and N2 is the selected / active cell, as it is after the edit
# is made interactively in Excel.
$ws.Range("N2").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 7
